$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Scottish bank-holiday dates (serial date numbers) appended below the
# existing list, continuing at row 113 through row 139.
$dates = @(
    42738, 42737, 42839, 42856, 42884, 42954, 43069, 43094, 43095, 43101,
    43102, 43189, 43227, 43248, 43318, 43434, 43459, 43460, 43466, 43467,
    43574, 43591, 43612, 43682, 43801, 43824, 43825
)

$startRow = 113
$endRow = $startRow + $dates.Length - 1

# Seed the new range with the same cell formatting (date number format) as
# the last existing data row (A112) before writing values, so the copied
# style is reused instead of a brand-new style being minted per cell.
$ws.Range("A112").Copy($ws.Range("A$startRow`:A$endRow"))

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
}

# Move the selection to the new last cell, matching the saved view state.
$ws.Range("A$endRow").Select()
